$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.277420333333333
$ws.Range("H2").Value = 3.832261
$ws.Range("I2").Value = 0.01913942624337554
$ws.Range("J2").Value = 0.01913942624337554
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 19.816421
$ws.Range("N2").Value = 59.449263
$ws.Range("O2").Value = 0.4265240049466206
$ws.Range("P2").Value = 0.4265240049466206
$ws.Range("Q2").Value = 25.31389911929367
$ws.Range("R2").Value = 227.825092073643
$ws.Range("S2").Value = 0.00816342473370499
$ws.Range("T2").Value = 0.008163424733704992

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.277420333333333
$ws.Range("H3").Value = 3.832261
$ws.Range("I3").Value = 0.01913942624337554
$ws.Range("J3").Value = 0.01913942624337554
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 18.17573333333333
$ws.Range("N3").Value = 54.5272
$ws.Range("O3").Value = 0.3912102278294916
$ws.Range("P3").Value = 0.3912102278294917
$ws.Range("Q3").Value = 23.21805133324445
$ws.Range("R3").Value = 208.9624619992
$ws.Range("S3").Value = 0.007487539301196697
$ws.Range("T3").Value = 0.007487539301196699

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.277420333333333
$ws.Range("H4").Value = 3.832261
$ws.Range("I4").Value = 0.01913942624337554
$ws.Range("J4").Value = 0.01913942624337554
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 6.429072666666666
$ws.Range("N4").Value = 19.287218
$ws.Range("O4").Value = 0.1383778545015528
$ws.Range("P4").Value = 0.1383778545015528
$ws.Range("Q4").Value = 8.212628148877556
$ws.Range("R4").Value = 73.91365333989799
$ws.Range("S4").Value = 0.002648472739949022
$ws.Range("T4").Value = 0.002648472739949023

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.277420333333333
$ws.Range("H5").Value = 3.832261
$ws.Range("I5").Value = 0.01913942624337554
$ws.Range("J5").Value = 0.01913942624337554
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.039044333333333
$ws.Range("N5").Value = 6.117133
$ws.Range("O5").Value = 0.04388791272233494
$ws.Range("P5").Value = 0.04388791272233494
$ws.Range("Q5").Value = 2.604716691968112
$ws.Range("R5").Value = 23.442450227713
$ws.Range("S5").Value = 0.0008399894685248326
$ws.Range("T5").Value = 0.0008399894685248327

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 45.44725166666667
$ws.Range("H6").Value = 136.341755
$ws.Range("I6").Value = 0.6809303864519871
$ws.Range("J6").Value = 0.6809303864519872
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 19.816421
$ws.Range("N6").Value = 59.449263
$ws.Range("O6").Value = 0.4265240049466206
$ws.Range("P6").Value = 0.4265240049466206
$ws.Range("Q6").Value = 900.6018723196185
$ws.Range("R6").Value = 8105.416850876566
$ws.Range("S6").Value = 0.2904331555193516
$ws.Range("T6").Value = 0.2904331555193517

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 45.44725166666667
$ws.Range("H7").Value = 136.341755
$ws.Range("I7").Value = 0.6809303864519871
$ws.Range("J7").Value = 0.6809303864519872
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 18.17573333333333
$ws.Range("N7").Value = 54.5272
$ws.Range("O7").Value = 0.3912102278294916
$ws.Range("P7").Value = 0.3912102278294917
$ws.Range("Q7").Value = 826.0371270262223
$ws.Range("R7").Value = 7434.334143236
$ws.Range("S7").Value = 0.2663869316199056
$ws.Range("T7").Value = 0.2663869316199058

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 45.44725166666667
$ws.Range("H8").Value = 136.341755
$ws.Range("I8").Value = 0.6809303864519871
$ws.Range("J8").Value = 0.6809303864519872
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 6.429072666666666
$ws.Range("N8").Value = 19.287218
$ws.Range("O8").Value = 0.1383778545015528
$ws.Range("P8").Value = 0.1383778545015528
$ws.Range("Q8").Value = 292.1836834652877
$ws.Range("R8").Value = 2629.65315118759
$ws.Range("S8").Value = 0.0942256859421392
$ws.Range("T8").Value = 0.09422568594213925

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 45.44725166666667
$ws.Range("H9").Value = 136.341755
$ws.Range("I9").Value = 0.6809303864519871
$ws.Range("J9").Value = 0.6809303864519872
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.039044333333333
$ws.Range("N9").Value = 6.117133
$ws.Range("O9").Value = 0.04388791272233494
$ws.Range("P9").Value = 0.04388791272233494
$ws.Range("Q9").Value = 92.66896097649057
$ws.Range("R9").Value = 834.020648788415
$ws.Range("S9").Value = 0.02988461337059061
$ws.Range("T9").Value = 0.02988461337059062

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.4966396666666666
$ws.Range("H10").Value = 1.489919
$ws.Range("I10").Value = 0.007441088905245192
$ws.Range("J10").Value = 0.007441088905245193
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 19.816421
$ws.Range("N10").Value = 59.449263
$ws.Range("O10").Value = 0.4265240049466206
$ws.Range("P10").Value = 0.4265240049466206
$ws.Range("Q10").Value = 9.841620719966334
$ws.Range("R10").Value = 88.574586479697
$ws.Range("S10").Value = 0.003173803041029044
$ws.Range("T10").Value = 0.003173803041029045

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.4966396666666666
$ws.Range("H11").Value = 1.489919
$ws.Range("I11").Value = 0.007441088905245192
$ws.Range("J11").Value = 0.007441088905245193
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 18.17573333333333
$ws.Range("N11").Value = 54.5272
$ws.Range("O11").Value = 0.3912102278294916
$ws.Range("P11").Value = 0.3912102278294917
$ws.Range("Q11").Value = 9.026790144088888
$ws.Range("R11").Value = 81.2411112968
$ws.Range("S11").Value = 0.002911030085920474
$ws.Range("T11").Value = 0.002911030085920475

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.4966396666666666
$ws.Range("H12").Value = 1.489919
$ws.Range("I12").Value = 0.007441088905245192
$ws.Range("J12").Value = 0.007441088905245193
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 6.429072666666666
$ws.Range("N12").Value = 19.287218
$ws.Range("O12").Value = 0.1383778545015528
$ws.Range("P12").Value = 0.1383778545015528
$ws.Range("Q12").Value = 3.192932506149111
$ws.Range("R12").Value = 28.736392555342
$ws.Range("S12").Value = 0.001029681917863138
$ws.Range("T12").Value = 0.001029681917863139

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.4966396666666666
$ws.Range("H13").Value = 1.489919
$ws.Range("I13").Value = 0.007441088905245192
$ws.Range("J13").Value = 0.007441088905245193
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 2.039044333333333
$ws.Range("N13").Value = 6.117133
$ws.Range("O13").Value = 0.04388791272233494
$ws.Range("P13").Value = 0.04388791272233494
$ws.Range("Q13").Value = 1.012670298025222
$ws.Range("R13").Value = 9.114032682227
$ws.Range("S13").Value = 0.0003265738604325358
$ws.Range("T13").Value = 0.0003265738604325359

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 19.52156333333333
$ws.Range("H14").Value = 58.56469
$ws.Range("I14").Value = 0.2924890983993922
$ws.Range("J14").Value = 0.2924890983993922
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 19.816421
$ws.Range("N14").Value = 59.449263
$ws.Range("O14").Value = 0.4265240049466206
$ws.Range("P14").Value = 0.4265240049466206
$ws.Range("Q14").Value = 386.8475175914967
$ws.Range("R14").Value = 3481.62765832347
$ws.Range("S14").Value = 0.124753621652535
$ws.Range("T14").Value = 0.124753621652535

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 19.52156333333333
$ws.Range("H15").Value = 58.56469
$ws.Range("I15").Value = 0.2924890983993922
$ws.Range("J15").Value = 0.2924890983993922
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 18.17573333333333
$ws.Range("N15").Value = 54.5272
$ws.Range("O15").Value = 0.3912102278294916
$ws.Range("P15").Value = 0.3912102278294917
$ws.Range("Q15").Value = 354.8187293964444
$ws.Range("R15").Value = 3193.368564568
$ws.Range("S15").Value = 0.1144247268224688
$ws.Range("T15").Value = 0.1144247268224689

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 19.52156333333333
$ws.Range("H16").Value = 58.56469
$ws.Range("I16").Value = 0.2924890983993922
$ws.Range("J16").Value = 0.2924890983993922
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 6.429072666666666
$ws.Range("N16").Value = 19.287218
$ws.Range("O16").Value = 0.1383778545015528
$ws.Range("P16").Value = 0.1383778545015528
$ws.Range("Q16").Value = 125.5055492369356
$ws.Range("R16").Value = 1129.54994313242
$ws.Range("S16").Value = 0.04047401390160146
$ws.Range("T16").Value = 0.04047401390160148

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 19.52156333333333
$ws.Range("H17").Value = 58.56469
$ws.Range("I17").Value = 0.2924890983993922
$ws.Range("J17").Value = 0.2924890983993922
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 2.039044333333333
$ws.Range("N17").Value = 6.117133
$ws.Range("O17").Value = 0.04388791272233494
$ws.Range("P17").Value = 0.04388791272233494
$ws.Range("Q17").Value = 39.80533309264111
$ws.Range("R17").Value = 358.24799783377
$ws.Range("S17").Value = 0.01283673602278696
$ws.Range("T17").Value = 0.01283673602278696
